$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24092
$ws.Range("D2").Value = 35156499
$ws.Range("C3").Value = 60804
$ws.Range("D3").Value = 89956419
$ws.Range("C4").Value = 20679
$ws.Range("D4").Value = 30842893
$ws.Range("C5").Value = 5454
$ws.Range("D5").Value = 8155097
$ws.Range("C6").Value = 1090
$ws.Range("D6").Value = 1632697
$ws.Range("C7").Value = 67
$ws.Range("D7").Value = 100500
$ws.Range("C10").Value = 25848
$ws.Range("D10").Value = 35682192
$ws.Range("C11").Value = 6307
$ws.Range("D11").Value = 9222147
$ws.Range("C12").Value = 17838
$ws.Range("D12").Value = 26365522
$ws.Range("C13").Value = 5580
$ws.Range("D13").Value = 8327171
$ws.Range("C14").Value = 1326
$ws.Range("D14").Value = 1981619
$ws.Range("C15").Value = 239
$ws.Range("D15").Value = 355766
$ws.Range("C17").Value = 6495
$ws.Range("D17").Value = 8782368
$ws.Range("C18").Value = 8766
$ws.Range("D18").Value = 12766330
$ws.Range("C19").Value = 21796
$ws.Range("D19").Value = 32240515
$ws.Range("C20").Value = 6943
$ws.Range("D20").Value = 10369461
$ws.Range("C21").Value = 1676
$ws.Range("D21").Value = 2508668
$ws.Range("C22").Value = 246
$ws.Range("D22").Value = 368624
$ws.Range("C24").Value = 7537
$ws.Range("D24").Value = 10297479
$ws.Range("C25").Value = 4972
$ws.Range("D25").Value = 7249883
$ws.Range("C26").Value = 15360
$ws.Range("D26").Value = 22715206
$ws.Range("C27").Value = 5240
$ws.Range("D27").Value = 7829969
$ws.Range("C28").Value = 1264
$ws.Range("D28").Value = 1894837
$ws.Range("C29").Value = 186
$ws.Range("D29").Value = 279000
$ws.Range("C31").Value = 5348
$ws.Range("D31").Value = 7174117
$ws.Range("C32").Value = 1756
$ws.Range("D32").Value = 2530919
$ws.Range("C33").Value = 4642
$ws.Range("D33").Value = 6815591
$ws.Range("C34").Value = 1887
$ws.Range("D34").Value = 2807851
$ws.Range("C35").Value = 495
$ws.Range("D35").Value = 739041
$ws.Range("C38").Value = 1195
$ws.Range("D38").Value = 1627872
$ws.Range("C39").Value = 11230
$ws.Range("D39").Value = 16367634
$ws.Range("C40").Value = 34685
$ws.Range("D40").Value = 51234896
$ws.Range("C41").Value = 12860
$ws.Range("D41").Value = 19204635
$ws.Range("C42").Value = 3575
$ws.Range("D42").Value = 5350598
$ws.Range("C43").Value = 627
$ws.Range("D43").Value = 937953
$ws.Range("C46").Value = 10541
$ws.Range("D46").Value = 14335991
$ws.Range("C47").Value = 1016
$ws.Range("D47").Value = 1470198
$ws.Range("C48").Value = 3808
$ws.Range("D48").Value = 5613489
$ws.Range("C49").Value = 1419
$ws.Range("D49").Value = 2121464
$ws.Range("C50").Value = 439
$ws.Range("D50").Value = 656000
$ws.Range("C52").Value = 2533
$ws.Range("D52").Value = 3521757
$ws.Range("C53").Value = 362
$ws.Range("D53").Value = 526784
$ws.Range("C54").Value = 970
$ws.Range("D54").Value = 1436744
$ws.Range("C55").Value = 397
$ws.Range("D55").Value = 593072
$ws.Range("C59").Value = 475
$ws.Range("D59").Value = 675625
$ws.Range("C60").Value = 10140
$ws.Range("D60").Value = 14728746
$ws.Range("C61").Value = 30932
$ws.Range("D61").Value = 45626014
$ws.Range("C62").Value = 10789
$ws.Range("D62").Value = 16119058
$ws.Range("C63").Value = 2993
$ws.Range("D63").Value = 4478160
$ws.Range("C64").Value = 531
$ws.Range("D64").Value = 795431
$ws.Range("C66").Value = 6
$ws.Range("D66").Value = 7787
$ws.Range("C67").Value = 10053
$ws.Range("D67").Value = 13461327
$ws.Range("C68").Value = 2774
$ws.Range("D68").Value = 4045452
$ws.Range("C69").Value = 7509
$ws.Range("D69").Value = 11073769
$ws.Range("C70").Value = 2660
$ws.Range("D70").Value = 3974483
$ws.Range("C71").Value = 878
$ws.Range("D71").Value = 1315049
$ws.Range("C72").Value = 175
$ws.Range("D72").Value = 261112
$ws.Range("C74").Value = 2950
$ws.Range("D74").Value = 4016747
$ws.Range("C75").Value = 886
$ws.Range("D75").Value = 1284391
$ws.Range("C76").Value = 3031
$ws.Range("D76").Value = 4476858
$ws.Range("C77").Value = 1207
$ws.Range("D77").Value = 1805939
$ws.Range("C78").Value = 418
$ws.Range("D78").Value = 627000
$ws.Range("C81").Value = 1810
$ws.Range("D81").Value = 2432059
$ws.Range("C82").Value = 311
$ws.Range("D82").Value = 462689
$ws.Range("C83").Value = 110
$ws.Range("D83").Value = 164610
$ws.Range("C85").Value = 17
$ws.Range("D85").Value = 25500
$ws.Range("C87").Value = 7128
$ws.Range("D87").Value = 10425054
$ws.Range("C88").Value = 20512
$ws.Range("D88").Value = 30341239
$ws.Range("C89").Value = 6756
$ws.Range("D89").Value = 10095322
$ws.Range("C90").Value = 1788
$ws.Range("D90").Value = 2676771
$ws.Range("C91").Value = 288
$ws.Range("D91").Value = 431810
$ws.Range("C92").Value = 24
$ws.Range("D92").Value = 36000
$ws.Range("C94").Value = 6403
$ws.Range("D94").Value = 8621161
$ws.Range("C95").Value = 19580
$ws.Range("D95").Value = 28434203
$ws.Range("C96").Value = 45564
$ws.Range("D96").Value = 67186825
$ws.Range("C97").Value = 14616
$ws.Range("D97").Value = 21815946
$ws.Range("C98").Value = 3926
$ws.Range("D98").Value = 5871908
$ws.Range("C99").Value = 663
$ws.Range("D99").Value = 992862
$ws.Range("C101").Value = 9
$ws.Range("D101").Value = 13500
$ws.Range("C102").Value = 16771
$ws.Range("D102").Value = 22775050
$ws.Range("C103").Value = 22396
$ws.Range("D103").Value = 32551921
$ws.Range("C104").Value = 50723
$ws.Range("D104").Value = 74715644
$ws.Range("C105").Value = 15914
$ws.Range("D105").Value = 23717248
$ws.Range("C106").Value = 4070
$ws.Range("D106").Value = 6079501
$ws.Range("C107").Value = 667
$ws.Range("D107").Value = 997554
$ws.Range("C108").Value = 31
$ws.Range("D108").Value = 44728
$ws.Range("C110").Value = 19904
$ws.Range("D110").Value = 26845223
$ws.Range("C111").Value = 8701
$ws.Range("D111").Value = 12706185
$ws.Range("C112").Value = 22577
$ws.Range("D112").Value = 33399155
$ws.Range("C113").Value = 7872
$ws.Range("D113").Value = 11749723
$ws.Range("C114").Value = 1914
$ws.Range("D114").Value = 2861601
$ws.Range("C115").Value = 273
$ws.Range("D115").Value = 406762
$ws.Range("C118").Value = 7130
$ws.Range("D118").Value = 9702412
$ws.Range("C119").Value = 21830
$ws.Range("D119").Value = 31733694
$ws.Range("C120").Value = 53844
$ws.Range("D120").Value = 79382968
$ws.Range("C121").Value = 16277
$ws.Range("D121").Value = 24301440
$ws.Range("C122").Value = 4048
$ws.Range("D122").Value = 6055497
$ws.Range("C123").Value = 816
$ws.Range("D123").Value = 1222212
$ws.Range("C126").Value = 18558
$ws.Range("D126").Value = 25486511
$ws.Range("C127").Value = 30060
$ws.Range("D127").Value = 44011111
$ws.Range("C128").Value = 90469
$ws.Range("D128").Value = 133983971
$ws.Range("C129").Value = 40169
$ws.Range("D129").Value = 60029201
$ws.Range("C130").Value = 12503
$ws.Range("D130").Value = 18726517
$ws.Range("C131").Value = 2576
$ws.Range("D131").Value = 3856006
$ws.Range("C132").Value = 136
$ws.Range("D132").Value = 203212
$ws.Range("C135").Value = 29560
$ws.Range("D135").Value = 41176828
